$wb = $excel.ActiveWorkbook

# Relabel shared strings: "light goods" -> "van", "heavy goods" -> "lorry"
# Applied identically on every worksheet since they share the same row labels in column A

$ws = $wb.Worksheets.Item("mean")
$ws.Range("A3").Value = "van"
$ws.Range("A4").Value = "lorry"

$ws = $wb.Worksheets.Item("median")
$ws.Range("A3").Value = "van"
$ws.Range("A4").Value = "lorry"

$ws = $wb.Worksheets.Item("lower 5")
$ws.Range("A3").Value = "van"
$ws.Range("A4").Value = "lorry"

$ws = $wb.Worksheets.Item("upper 95")
$ws.Range("A3").Value = "van"
$ws.Range("A4").Value = "lorry"

# Updated numeric results (re-run of the underlying rate calculation)

$ws = $wb.Worksheets.Item("mean")
$ws.Range("B2").Value = 2.967330591485387
$ws.Range("C2").Value = 1.572154356438868
$ws.Range("D2").Value = 2.331115657530493
$ws.Range("E2").Value = 1.3276926034946215
$ws.Range("F2").Value = 2.0077068422803177
$ws.Range("B3").Value = 1.8884764605965445
$ws.Range("C3").Value = 1.073521283726366
$ws.Range("D3").Value = 1.5430219092042556
$ws.Range("E3").Value = 0.8696898463479509
$ws.Range("F3").Value = 1.3322570700669945
$ws.Range("B4").Value = 4.489380915224699
$ws.Range("C4").Value = 3.9658562595259026
$ws.Range("D4").Value = 6.885556875025202
$ws.Range("E4").Value = 4.715036952981986
$ws.Range("F4").Value = 4.65299756713466
$ws.Range("B5").Value = 1.530153226496554
$ws.Range("C5").Value = 0.9099386208418274
$ws.Range("D5").Value = 0.7617813104174076
$ws.Range("E5").Value = 0.2619264799458266
$ws.Range("F5").Value = 0.6375000049661662
$ws.Range("B6").Value = 18.91264674909875
$ws.Range("C6").Value = 15.892451230399145
$ws.Range("D6").Value = 14.110878657838642
$ws.Range("E6").Value = 11.585507983853423
$ws.Range("F6").Value = 14.456107026260627
$ws.Range("B7").Value = 5.77432994056328
$ws.Range("C7").Value = 0.5964434916318531
$ws.Range("D7").Value = 0.7990041624664719
$ws.Range("E7").Value = 0.2736051279665286
$ws.Range("F7").Value = 0.48660212009861353

$ws = $wb.Worksheets.Item("median")
$ws.Range("B2").Value = 2.9688824908236
$ws.Range("C2").Value = 1.5751364784404562
$ws.Range("D2").Value = 2.325591326820291
$ws.Range("E2").Value = 1.3257113673321643
$ws.Range("F2").Value = 2.007834583155332
$ws.Range("B3").Value = 1.8583336564066275
$ws.Range("C3").Value = 1.0416938770681399
$ws.Range("D3").Value = 1.51394803778331
$ws.Range("E3").Value = 0.8534775964112733
$ws.Range("F3").Value = 1.3079806389195228
$ws.Range("B4").Value = 4.3711678787806125
$ws.Range("C4").Value = 3.888904226441387
$ws.Range("D4").Value = 6.742130616565655
$ws.Range("E4").Value = 4.622401542894465
$ws.Range("F4").Value = 4.564394847458054
$ws.Range("B5").Value = 0.9508879781216066
$ws.Range("C5").Value = 0.5650544840536169
$ws.Range("D5").Value = 0.4786222382546181
$ws.Range("E5").Value = 0.16345406325793554
$ws.Range("F5").Value = 0.3960897116547871
$ws.Range("B6").Value = 18.553754299725135
$ws.Range("C6").Value = 15.600868472710701
$ws.Range("D6").Value = 13.79691699714919
$ws.Range("E6").Value = 11.409630988202363
$ws.Range("F6").Value = 14.240512272954453
$ws.Range("B7").Value = 4.788124047237519
$ws.Range("C7").Value = 0.5187953539396203
$ws.Range("D7").Value = 0.6906228565394069
$ws.Range("E7").Value = 0.2366347187545119
$ws.Range("F7").Value = 0.4342338821864423

$ws = $wb.Worksheets.Item("lower 5")
$ws.Range("B2").Value = 2.8065431985808234
$ws.Range("C2").Value = 1.4522533401944457
$ws.Range("D2").Value = 2.1645677544684423
$ws.Range("E2").Value = 1.235478499149242
$ws.Range("F2").Value = 1.9400439699452396
$ws.Range("B3").Value = 1.3956634962249779
$ws.Range("C3").Value = 0.789513207206368
$ws.Range("D3").Value = 1.128458966099467
$ws.Range("E3").Value = 0.6348784545675036
$ws.Range("F3").Value = 0.9993417248253278
$ws.Range("B4").Value = 3.1549331294461815
$ws.Range("C4").Value = 2.7084125601633673
$ws.Range("D4").Value = 4.660763384434445
$ws.Range("E4").Value = 3.1445295569379432
$ws.Range("F4").Value = 3.282112394512766
$ws.Range("B5").Value = 0.19896259211327433
$ws.Range("C5").Value = 0.12227519789536709
$ws.Range("D5").Value = 0.10232724740560545
$ws.Range("E5").Value = 0.03381533253926477
$ws.Range("F5").Value = 0.08553988996918686
$ws.Range("B6").Value = 14.221861010692573
$ws.Range("C6").Value = 12.13455804716933
$ws.Range("D6").Value = 10.17562608888569
$ws.Range("E6").Value = 8.779454991823009
$ws.Range("F6").Value = 11.255807640233702
$ws.Range("B7").Value = 1.8545870468950356
$ws.Range("C7").Value = 0.1951394520213546
$ws.Range("D7").Value = 0.26737655517164854
$ws.Range("E7").Value = 0.09669112461449791
$ws.Range("F7").Value = 0.191243632229466

$ws = $wb.Worksheets.Item("upper 95")
$ws.Range("B2").Value = 3.1299293353803437
$ws.Range("C2").Value = 1.6950859832955394
$ws.Range("D2").Value = 2.5084591676743
$ws.Range("E2").Value = 1.4236650790520728
$ws.Range("F2").Value = 2.0774558933952383
$ws.Range("B3").Value = 2.456842359046504
$ws.Range("C3").Value = 1.4596836803751954
$ws.Range("D3").Value = 2.0670007277378395
$ws.Range("E3").Value = 1.1603756792223459
$ws.Range("F3").Value = 1.722506138432508
$ws.Range("B4").Value = 6.075392802512143
$ws.Range("C4").Value = 5.511009453625714
$ws.Range("D4").Value = 9.629166738096071
$ws.Range("E4").Value = 6.612648167045657
$ws.Range("F4").Value = 6.251193697746722
$ws.Range("B5").Value = 4.734756045397359
$ws.Range("C5").Value = 2.8005688133817963
$ws.Range("D5").Value = 2.2555708341947636
$ws.Range("E5").Value = 0.7865798184454649
$ws.Range("F5").Value = 1.9212283514551105
$ws.Range("B6").Value = 24.30005028998928
$ws.Range("C6").Value = 20.676067237783922
$ws.Range("D6").Value = 19.028636595909386
$ws.Range("E6").Value = 14.898712625482954
$ws.Range("F6").Value = 18.08274608283083
$ws.Range("B7").Value = 13.059881778385057
$ws.Range("C7").Value = 1.2828585390225573
$ws.Range("D7").Value = 1.6650479463189323
$ws.Range("E7").Value = 0.6106254097354231
$ws.Range("F7").Value = 1.0140680683211465
